{"js": "// Update the 5x5 \"two-digit \u00f7 one-digit\" practice table: each non-blank\n// cell's expression is replaced with a new one (same \"a\u00f7b=\" shape), per\n// the authoritative before -> after text map below. Cells are matched by\n// their *current* text, so duplicate expressions (e.g. \"66\u00f74=\" which\n// appears twice, each mapping to a different replacement) are resolved\n// positionally in document order, exactly like the source edit.\nconst replacements = [\n  [\"22\u00f75=\", \"39\u00f75=\"],\n  [\"96\u00f76=\", \"72\u00f75=\"],\n  [\"50\u00f74=\", \"26\u00f79=\"],\n  [\"99\u00f79=\", \"98\u00f76=\"],\n  [\"50\u00f78=\", \"46\u00f78=\"],\n  [\"26\u00f75=\", \"13\u00f77=\"],\n  [\"81\u00f77=\", \"67\u00f74=\"],\n  [\"75\u00f75=\", \"16\u00f74=\"],\n  [\"35\u00f77=\", \"51\u00f77=\"],\n  [\"13\u00f76=\", \"20\u00f77=\"],\n  [\"66\u00f74=\", \"15\u00f76=\"],\n  [\"73\u00f78=\", \"73\u00f73=\"],\n  [\"15\u00f79=\", \"94\u00f77=\"],\n  [\"84\u00f74=\", \"30\u00f72=\"],\n  [\"17\u00f78=\", \"23\u00f78=\"],\n  [\"15\u00f73=\", \"78\u00f76=\"],\n  [\"94\u00f74=\", \"43\u00f73=\"],\n  [\"54\u00f76=\", \"65\u00f78=\"],\n  [\"16\u00f74=\", \"64\u00f79=\"],\n  [\"37\u00f74=\", \"24\u00f75=\"],\n  [\"66\u00f74=\", \"14\u00f73=\"],\n  [\"41\u00f73=\", \"75\u00f77=\"],\n  [\"10\u00f79=\", \"47\u00f76=\"],\n  [\"48\u00f78=\", \"36\u00f74=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nlet next = 0;\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    const current = grid[r][c];\n    if (next < replacements.length && current === replacements[next][0]) {\n      grid[r][c] = replacements[next][1];\n      next++;\n    }\n  }\n}\n\ntable.values = grid;\nawait context.sync();\n", "ps1": "# Update the 5x5 \"two-digit \u00f7 one-digit\" practice table: each non-blank\n# cell's expression is replaced with a new one (same \"a\u00f7b=\" shape), per\n# the authoritative before -> after text map below. Cells are matched by\n# their *current* text, so duplicate expressions (e.g. \"66\u00f74=\" which\n# appears twice, each mapping to a different replacement) are resolved\n# positionally in document (row-major) order, exactly like the source edit.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$map = @(\n  @(\"22\u00f75=\", \"39\u00f75=\"),\n  @(\"96\u00f76=\", \"72\u00f75=\"),\n  @(\"50\u00f74=\", \"26\u00f79=\"),\n  @(\"99\u00f79=\", \"98\u00f76=\"),\n  @(\"50\u00f78=\", \"46\u00f78=\"),\n  @(\"26\u00f75=\", \"13\u00f77=\"),\n  @(\"81\u00f77=\", \"67\u00f74=\"),\n  @(\"75\u00f75=\", \"16\u00f74=\"),\n  @(\"35\u00f77=\", \"51\u00f77=\"),\n  @(\"13\u00f76=\", \"20\u00f77=\"),\n  @(\"66\u00f74=\", \"15\u00f76=\"),\n  @(\"73\u00f78=\", \"73\u00f73=\"),\n  @(\"15\u00f79=\", \"94\u00f77=\"),\n  @(\"84\u00f74=\", \"30\u00f72=\"),\n  @(\"17\u00f78=\", \"23\u00f78=\"),\n  @(\"15\u00f73=\", \"78\u00f76=\"),\n  @(\"94\u00f74=\", \"43\u00f73=\"),\n  @(\"54\u00f76=\", \"65\u00f78=\"),\n  @(\"16\u00f74=\", \"64\u00f79=\"),\n  @(\"37\u00f74=\", \"24\u00f75=\"),\n  @(\"66\u00f74=\", \"14\u00f73=\"),\n  @(\"41\u00f73=\", \"75\u00f77=\"),\n  @(\"10\u00f79=\", \"47\u00f76=\"),\n  @(\"48\u00f78=\", \"36\u00f74=\")\n)\n\n$idx = 0\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cur = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($idx -lt $map.Length -and $cur -eq $map[$idx][0]) {\n      $cell.Range.Text = $map[$idx][1]\n      $idx++\n    }\n  }\n}\n"}
